# "split snapshot and release for auto deploying"
#
# The automated test-report writer re-wrote column F (the per-step
# Résultat column) of the single worksheet with the outcome of each
# Jenkins scenario: a green "Succès" for the two steps that passed and a
# red "Échec : ..." message (with the concrete failure reason) for the
# three steps that failed. Reproduce that here using the Excel object
# model: set the cell value/text and colour the font accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ColorIndex 10 resolves to RGB 008000 (the same green already used by
# the workbook's pre-existing "success" font), ColorIndex 3 resolves to
# RGB FF0000 (the same red already used by the "failure" font).
$greenColorIndex = 10
$redColorIndex   = 3

# Row 2 -> "Ouverture DEMO" step: succeeded
$ws.Range("F2").Value = "Succès"
$ws.Range("F2").Font.ColorIndex = $greenColorIndex

# Row 3 -> "The city is Paris!!" step: failed
$ws.Range("F3").Value = "Échec : The city is Paris!!"
$ws.Range("F3").Font.ColorIndex = $redColorIndex

# Row 4 -> "Saisie « Input Select field » dans demo" step: failed
$ws.Range("F4").Value = "Échec : Saisie « Input Select field » dans demo."
$ws.Range("F4").Font.ColorIndex = $redColorIndex

# Row 5 -> second "Ouverture DEMO"-style step: succeeded
$ws.Range("F5").Value = "Succès"
$ws.Range("F5").Font.ColorIndex = $greenColorIndex

# Row 6 -> "Accès à l'action « no exist element » dans demo" step: failed
$ws.Range("F6").Value = "Échec : Accès à l'action « no exist element » dans demo."
$ws.Range("F6").Font.ColorIndex = $redColorIndex
